# Apply the "preparation_temperature" -> "preparation_condition" and
# "storage_temperature" -> "storage_method" rework to the sample-block
# workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Rename the two lookup-list worksheets (sheetId / position unchanged)
# ---------------------------------------------------------------------
$wsPrepCond = $wb.Worksheets.Item("preparation_temperature list")
$wsPrepCond.Name = "preparation_condition list"

$wsStorMethod = $wb.Worksheets.Item("storage_temperature list")
$wsStorMethod.Name = "storage_method list"

# ---------------------------------------------------------------------
# 2. Replace the values in the "preparation_condition list" sheet
#    (7 values, was 8)
# ---------------------------------------------------------------------
$prepValues = @(
    "frozen in liquid nitrogen",
    "frozen in liquid nitrogen vapor",
    "frozen in ice",
    "frozen in dry ice",
    "frozen at -20 C",
    "ambient temperature",
    "unknown"
)
for ($i = 0; $i -lt $prepValues.Length; $i++) {
    $wsPrepCond.Cells.Item($i + 1, 1).Value = $prepValues[$i]
}
# remove the old 8th row so the list is exactly A1:A7
$wsPrepCond.Rows.Item(8).Delete()

# ---------------------------------------------------------------------
# 3. Replace the values in the "storage_method list" sheet
#    (11 values, was 12)
# ---------------------------------------------------------------------
$storValues = @(
    "frozen in liquid nitrogen",
    "frozen in liquid nitrogen vapor",
    "frozen in ice",
    "frozen in dry ice",
    "frozen at -80 C",
    "frozen at -20 C",
    "refrigerator",
    "ambient temperature",
    "incubated at 37 C",
    "none",
    "unknown"
)
for ($i = 0; $i -lt $storValues.Length; $i++) {
    $wsStorMethod.Cells.Item($i + 1, 1).Value = $storValues[$i]
}
# remove the old 12th row so the list is exactly A1:A11
$wsStorMethod.Rows.Item(12).Delete()

# ---------------------------------------------------------------------
# 4. Update the main data sheet: header text, comments and validations
# ---------------------------------------------------------------------
$wsMain = $wb.Worksheets.Item("Export as TSV")

# Column M header ("preparation_temperature" -> "preparation_condition")
$wsMain.Cells.Item(1, 13).Value = "preparation_condition"

# Column Q header ("storage_temperature" -> "storage_method")
$wsMain.Cells.Item(1, 17).Value = "storage_method"

# Comment on M1
[void]$wsMain.Cells.Item(1, 13).Comment.Text("The condition under which the preparation occurred, such as whether the sample was placed in dry ice during the preparation.")

# Comment on Q1
[void]$wsMain.Cells.Item(1, 17).Comment.Text("The method by which the sample was stored, after preparation and before the assay was performed.")

# Data validation for column M (preparation_condition)
$dvM = $wsMain.Range("M2:M1048576").Validation
$dvM.Modify(3, 1, 1, "'preparation_condition list'!`$A`$1:`$A`$7")
$dvM.ErrorTitle = "Value must come from list"
$dvM.ErrorMessage = "Value must come from preparation_condition list."

# Data validation for column Q (storage_method)
$dvQ = $wsMain.Range("Q2:Q1048576").Validation
$dvQ.Modify(3, 1, 1, "'storage_method list'!`$A`$1:`$A`$11")
$dvQ.ErrorTitle = "Value must come from list"
$dvQ.ErrorMessage = "Value must come from storage_method list."
